$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gesamtergebnis")

$ws.Range("C3").Value = 21.19
$ws.Range("D3").Value = 1.88

$ws.Range("D9").Select()
